$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "EpubBooks"

# --- Add the new sheet right after EpubBooks ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "People to test"

# Match the authored page margins (0.75in/1in/0.5in -> 54/72/36 points)
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# --- Populate the new sheet with the people data ---
$ws2.Range("A2").Value = "David"
$ws2.Range("A3").Value = "William"
$ws2.Range("A4").Value = "Jai"

# Column A width -> stored xlsx width of 15.5 (COM ColumnWidth differs from
# the stored character width by the default font padding, 5/6 chars)
$ws2.Columns.Item(1).ColumnWidth = 44/3

# Make the new sheet the active one, with the selection sitting just below
# the data, matching the authored workbook state.
$ws2.Activate()
$ws2.Range("A5").Select() | Out-Null
